$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sat Dec 16 23_03_43 2023"
$ws.Range("B2").Value = "loc"
$ws.Range("C2").Value = -1
